$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append the "Twenty Questions" answer key: animal name + its row index.
# Rows 18-20 are intentionally left blank (matches the source data layout),
# new data starts at row 21.
# ---------------------------------------------------------------------------
$animals = @(
  "Dog", "Human", "Pig", "Bat", "Tiger", "Rat", "Deer", "Giraffe",
  "Monkey", "Chicken", "Turtle", "Lizard", "Shark", "Penguin", "Bird", "Snake"
)

$startRow = 21
for ($i = 0; $i -lt $animals.Length; $i++) {
  $r = $startRow + $i
  $ws.Range("A$r").Value = $i
  $ws.Range("B$r").Value = $animals[$i]
}

# ---------------------------------------------------------------------------
# Widen columns C:G for readability (values chosen to land on the closest
# width the host engine can represent).
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 15.3333
$ws.Columns.Item(4).ColumnWidth = 15.3333
$ws.Columns.Item(5).ColumnWidth = 16.1667
$ws.Columns.Item(6).ColumnWidth = 19.0
$ws.Columns.Item(7).ColumnWidth = 16.8333

# ---------------------------------------------------------------------------
# Leave the selection where the author last left off while finishing the
# report.
# ---------------------------------------------------------------------------
$ws.Range("C33").Select()
